# Slide 12, shape 1 ("object 2" - the title placeholder holding the
# "Query q = em.createQuery(...)" code line) needs:
#   1. Its text split so the single run ending the paragraph
#      ("Registration");) becomes three runs: "Registration", " r", "");"
#   2. Its position/size updated to the new xfrm values.
#
# Note: the shape has <a:spAutoFit/>, so changing its text can cause the
# engine to recompute Height automatically; set Left/Top/Width/Height
# *after* the text edits so the explicit values win.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(12)
$sh = $s.Shapes.Item(1)

$tr = $sh.TextFrame.TextRange

# The paragraph's last run currently reads: Registration");
# Split it into: "Registration" | " r" | "”);"
$lastRun = $tr.Runs(13, 1)
$lastRun.Text = "Registration"
$lastRun.InsertAfter(" r")

$midRun = $tr.Runs(14, 1)
$midRun.InsertAfter([char]8221 + ");")

# Reposition / resize the shape (values taken from the target xfrm, chosen
# so that the points-to-EMU round trip lands on the exact target EMUs).
$sh.Left   = 89.5
$sh.Top    = 415.25
$sh.Width  = 1443.4500732421875
$sh.Height = 48.9
